# core/static/fatality_data_for_upload.xlsx
#
# - Add new valid species "Pogoniulus atroflavus" (alphabetically just
#   before the existing "Pogoniulus bilineatus" entry) on the
#   "Valid Species" sheet.
# - Remove the five obsolete "Turdus philomelos*" taxon rows from the
#   same sheet (species re-classified / renamed upstream).
# - Keep the "Main" sheet's genus/species data-validation list in sync
#   with the new extent of the "Valid Species" column.

$wb = $excel.ActiveWorkbook
$validSheet = $wb.Worksheets.Item("Valid Species")
$mainSheet = $wb.Worksheets.Item("Main")

# The sheet is protected (read-only cell entry) - temporarily lift that
# so the list can be edited, restoring protection once done.
$validSheet.Unprotect()

# Insert "Pogoniulus atroflavus" right before "Pogoniulus bilineatus",
# keeping the list alphabetically sorted.
$validSheet.Rows.Item(3160).Insert()
$validSheet.Cells.Item(3160, 1).Value2 = "Pogoniulus atroflavus"

# Drop the five "Turdus philomelos ..." rows (now shifted down by one
# row because of the insert above).
$validSheet.Range("A4137:A4141").EntireRow.Delete()

$validSheet.Protect()

# Point the "Main" sheet genus/species validation list at the new last
# row of the (now shorter) Valid Species column.
$lastRow = $validSheet.Cells.Item($validSheet.Rows.Count, 1).End(-4162).Row
$dv = $mainSheet.Range("A2:A1048576").Validation
$dv.Formula1 = "='Valid Species'!A1:A" + $lastRow
